$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (rows 2-15), matching the updated export
$data = @(
    @("BN-BAG-17-choko-Photoroom", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\BN-BAG-17-choko\5_image_BN-BAG-17-choko-Photoroom.png", 510, 500, "min side 500 < 501"),
    @("BN-BAG-22-choko", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\BN-BAG-22-choko\005_image_BN-BAG-22-choko.jpg", 500, 604, "min side 500 < 501"),
    @("BN-BAG-28-choko", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\BN-BAG-28-choko\002_image_BN-BAG-28-choko.jpg", 500, 721, "min side 500 < 501"),
    @("BN-BAG-3-choko", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\BN-BAG-3-choko\002_image_BN-BAG-3-choko.jpg", 500, 602, "min side 500 < 501"),
    @("LC69055-76", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69055-76\5_image_LC69055-76.jpg", 977, 500, "min side 500 < 501"),
    @("LC69802-01", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-01\5_image_LC69802-01.jpg", 1154, 500, "min side 500 < 501"),
    @("LC69802-12", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-12\009_image_LC69802-12.jpg", 795, 500, "min side 500 < 501"),
    @("LC69802-12", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-12\4_image_LC69802-12.jpg", 866, 500, "min side 500 < 501"),
    @("LC69802-12", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-12\5_image_LC69802-12.jpg", 1177, 500, "min side 500 < 501"),
    @("LC69802-25", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-25\001_image_LC69802-25.jpg", 893, 500, "min side 500 < 501"),
    @("LC69802-25", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69802-25\4_image_LC69802-25.jpg", 897, 500, "min side 500 < 501"),
    @("LC69855-2058", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\LC69855-2058\001_image_LC69855-2058.jpg", 857, 500, "min side 500 < 501"),
    @("RCh-1313-4lx", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\RCh-1313-4lx\001_image_RCh-1313-4lx.jpg", 1183, 500, "min side 500 < 501"),
    @("TW-Lina-black-flo", "C:/Users/Asus/Desktop/1 20.01.2026/339 фото/2 Новая папка обрізані\TW-Lina-black-flo\7_image_TW-Lina-black-flo.jpg", 517, 500, "min side 500 < 501")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
